$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.690.83"
$ws.Range("E2").Value = "  +1.63%  "

$ws.Range("D3").Value = "1.637.33"
$ws.Range("E3").Value = "  +1.19%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.66"
$ws.Range("E5").Value = "  +0.80%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.507"
$ws.Range("E6").Value = "  +4.32%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("E8").Value = "  +2.78%  "

$ws.Range("E9").Value = "  +1.57%  "

$ws.Range("E10").Value = "  +3.02%  "

$ws.Range("E11").Value = "  +3.52%  "

$ws.Range("D12").Value = "1.865.13"
$ws.Range("E12").Value = "  +1.18%  "

$ws.Range("D13").Value = "1.631.93"
$ws.Range("E13").Value = "  +0.79%  "

$ws.Range("E14").Value = "  +2.70%  "

$ws.Range("E15").Value = "  +2.25%  "

$ws.Range("D16").Value = "26.687.32"
$ws.Range("E16").Value = "  +1.58%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.63"
$ws.Range("E17").Value = "  +2.08%  "

$ws.Range("E18").Value = "  +2.69%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "220.39"
$ws.Range("E19").Value = "  +9.54%  "

$ws.Range("E20").Value = "  +0.00%  "

$ws.Range("E21").Value = "  +1.25%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.47"
$ws.Range("E22").Value = "  +1.75%  "

$ws.Range("E23").Value = "  +3.03%  "

$ws.Range("E24").Value = "  +1.70%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.32"
$ws.Range("E25").Value = "  +3.02%  "

$ws.Range("E26").Value = "  -0.09%  "

$ws.Range("E27").Value = "  +1.66%  "

$ws.Range("E28").Value = "  +6.40%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.55"
$ws.Range("E29").Value = "  +2.65%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0513"
$ws.Range("E30").Value = "  -0.71%  "

$ws.Range("E31").Value = "  +0.14%  "

$ws.Range("E32").Value = "  +5.29%  "

$ws.Range("E33").Value = "  +2.81%  "

$ws.Range("E34").Value = "  +2.06%  "

$ws.Range("E35").Value = "  -0.65%  "

$ws.Range("D36").Value = "1.220.71"
$ws.Range("E36").Value = "  +3.64%  "

$ws.Range("E37").Value = "  +6.11%  "

$ws.Range("E38").Value = "  +1.32%  "

$ws.Range("E39").Value = "  +0.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.508"
$ws.Range("E40").Value = "  +2.85%  "

$ws.Range("E41").Value = "  -1.23%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.44"
$ws.Range("E42").Value = "  +2.12%  "

$ws.Range("E43").Value = "  +0.41%  "

$ws.Range("D44").Value = "1.774.66"
$ws.Range("E44").Value = "  +1.06%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "93.67"
$ws.Range("E45").Value = "  +0.97%  "

$ws.Range("E46").Value = "  +2.79%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.02"
$ws.Range("E47").Value = "  +2.54%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0514"
$ws.Range("E48").Value = "  +1.08%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.69"
$ws.Range("E49").Value = "  +5.91%  "

$ws.Range("E50").Value = "  +0.61%  "

$ws.Range("E51").Value = "  -0.01%  "
